$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These Price cells get new values that look like plain numbers. Force them
# to Text format first so Excel stores the exact locale-formatted digit string
# (e.g. "185.39") instead of silently converting it to a Number cell.
$textCells = @("D5","D6","D10","D14","D18","D20","D21","D22","D23","D28","D31","D32","D36","D37","D38","D39","D40","D41","D43","D44","D45","D46","D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "66.945.33"
$ws.Range("E2").Value = "  -1.32%  "
$ws.Range("D3").Value = "3.301.69"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "185.39"
$ws.Range("E5").Value = "  +0.99%  "
$ws.Range("D6").Value = "574.84"
$ws.Range("E6").Value = "  -0.95%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("D10").Value = "6.67"
$ws.Range("E10").Value = "  +1.33%  "
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").Value = "3.876.61"
$ws.Range("E12").Value = "  +1.13%  "
$ws.Range("E13").Value = "  -0.54%  "
$ws.Range("D14").Value = "27.38"
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("D15").Value = "67.230.75"
$ws.Range("E15").Value = "  -0.88%  "
$ws.Range("E16").Value = "  -0.91%  "
$ws.Range("D17").Value = "3.311.28"
$ws.Range("E17").Value = "  +1.35%  "
$ws.Range("D18").Value = "442.90"
$ws.Range("E18").Value = "  +10.04%  "
$ws.Range("E19").Value = "  +0.17%  "
$ws.Range("D20").Value = "13.51"
$ws.Range("E20").Value = "  +0.87%  "
$ws.Range("D21").Value = "7.75"
$ws.Range("E21").Value = "  +2.75%  "
$ws.Range("D22").Value = "74.23"
$ws.Range("E22").Value = "  +4.46%  "
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("D25").Value = "3.437.70"
$ws.Range("E25").Value = "  +0.88%  "
$ws.Range("E26").Value = "  +1.43%  "
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("D28").Value = "9.08"
$ws.Range("E28").Value = "  -3.97%  "
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("E30").Value = "  +1.39%  "
$ws.Range("D31").Value = "22.82"
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("D32").Value = "5.32"
$ws.Range("E32").Value = "  -2.57%  "
$ws.Range("E34").Value = "  -1.01%  "
$ws.Range("E35").Value = "  -1.69%  "
$ws.Range("D36").Value = "1.52"
$ws.Range("E36").Value = "  +4.99%  "
$ws.Range("D37").Value = "163.05"
$ws.Range("E37").Value = "  -0.85%  "
$ws.Range("D38").Value = "27.42"
$ws.Range("E38").Value = "  +1.12%  "
$ws.Range("D39").Value = "1.85"
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").Value = "0.786"
$ws.Range("E40").Value = "  -1.63%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").Value = "4.46"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").Value = "2.732.36"
$ws.Range("E42").Value = "  +1.93%  "
$ws.Range("D43").Value = "6.31"
$ws.Range("E43").Value = "  -0.45%  "
$ws.Range("D44").Value = "40.27"
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("D45").Value = "0.0672"
$ws.Range("E45").Value = "  -0.81%  "
$ws.Range("D46").Value = "24.72"
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("D48").Value = "327.79"
$ws.Range("E48").Value = "  -2.13%  "
$ws.Range("E49").Value = "  -0.15%  "
$ws.Range("E50").Value = "  +2.86%  "
$ws.Range("E51").Value = "  -0.86%  "
